$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.031.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4643"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3727"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("E9").Value = "  -2.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8654"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("E11").Value = "  -2.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07822"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.840.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.358"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.559"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008869"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.067.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.168"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.077.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.844"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.090"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.131"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.79%  "
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08866"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.967"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7313"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.134"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.499"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.25%  "
$ws.Range("E37").Value = "  -1.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01947"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05239"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.349"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.60%  "
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5171"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8567"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -15.15%  "
$ws.Range("E45").Value = "  -3.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4829"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.78%  "
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "102.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.17%  "
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06249"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.73%  "
